# Add "Set size" and "Is timeseries" columns to the PHENOTYPES and
# PHENOTYPES_EXAMPLE tables (trials-data template), and populate the
# example data for the two new columns on PHENOTYPES_EXAMPLE.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# PHENOTYPES_EXAMPLE sheet (table "Table13") - add columns + sample data
# ---------------------------------------------------------------------
$wsExample = $wb.Worksheets.Item("PHENOTYPES_EXAMPLE")
$loExample = $wsExample.ListObjects.Item(1)

$colExampleSetSize = $loExample.ListColumns.Add()
$wsExample.Range("K1").Value = "Set size"

$colExampleTimeseries = $loExample.ListColumns.Add()
$wsExample.Range("L1").Value = "Is timeseries"

$wsExample.Range("K2").Value = 3
$wsExample.Range("L2").Value = $false

$wsExample.Range("L3").Value = $true

$wsExample.Range("K4").Value = 3
$wsExample.Range("L4").Value = $false

$wsExample.Range("K5").Value = 1
$wsExample.Range("L5").Value = $false

# Scroll the view so column C is the left-most visible column, then
# select the new header cells (mirrors the manual edit that introduced
# these columns).
$wsExample.Activate()
$winExample = $excel.ActiveWindow
$winExample.ScrollColumn = 3
$wsExample.Range("K1:L1").Select()

# ---------------------------------------------------------------------
# PHENOTYPES sheet (table "Table136") - add the matching template columns
# ---------------------------------------------------------------------
$wsPhen = $wb.Worksheets.Item("PHENOTYPES")
$loPhen = $wsPhen.ListObjects.Item(1)

$colPhenSetSize = $loPhen.ListColumns.Add()
$wsPhen.Range("K1").Value = "Set size"

$colPhenTimeseries = $loPhen.ListColumns.Add()
$wsPhen.Range("L1").Value = "Is timeseries"

# Match the column widths Excel computed for the new columns.
$wsPhen.Columns.Item(11).ColumnWidth = 9.1
$wsPhen.Columns.Item(12).ColumnWidth = 13.6

$wsPhen.Range("K1:L1").Select()

# ---------------------------------------------------------------------
# Restore METADATA as the active sheet/tab (it was active before/after
# this edit and its view state is untouched by this change).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("METADATA").Activate()
